$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Row 11 ("Marking"): Right count 5 -> 4, Wrong count -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 ("Total"): total score 90 -> 72, and the "x / y" label updated accordingly
$ws.Range("B12").Value = 72
$ws.Range("E12").Value = "72 / 112"
